# ECEN361 Project Definition - "cleaned up headers on first page"
#
# 1. Replace the "<Group Name>" placeholder with the real group name ("STC Rover"),
#    and replace "<Team MEmber Names>" with the real member names ("Hannah, Seth, Chris"),
#    keeping the manual line break between them in its own run (matching how Word
#    itself splits a run when new text is typed after a <w:br/>).
# 2. Remove the instructional "Describe the key elements..." paragraph that used to
#    follow the "Lessons to Learn" heading.
# 3. Move the <w:lastRenderedPageBreak/> hint from the "Schedule" heading run to the
#    "Derived Requirements" heading run (the page break now falls there instead).

$d = $word.ActiveDocument

# --- helper: find the paragraph whose trimmed text equals $needle ----------------
function Get-ParaByText($doc, [string]$needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text -replace "[\r\v]", ""
        if ($t -eq $needle) {
            return $p
        }
    }
    return $null
}

# --- 1. Group name / team member names on the title page ------------------------
$subtitlePara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text -replace "[\r\v]", ""
    if ($t -eq "<Group Name><Team MEmber Names>") {
        $subtitlePara = $p
        break
    }
}

if ($subtitlePara -ne $null) {
    $xml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="308D359C" w14:textId="710691B3" w:rsidR="00856EB6" w:rsidRPr="00856EB6" w:rsidRDefault="00DE5B15" w:rsidP="00856EB6">
<w:pPr><w:pStyle w:val="Subtitle"/><w:jc w:val="right"/></w:pPr>
<w:r><w:t>STC Rover</w:t></w:r>
<w:r><w:br/></w:r>
<w:r><w:t>Hannah, Seth, Chris</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
    $subtitlePara.Range.InsertXML($xml)
}

# --- 2. Drop the "Describe the key elements..." instructional paragraph ---------
$instrPara = Get-ParaByText $d "Describe the key elements you’ll need to learn more about to complete your project. Include information such as “the team will learn about wireless and cloud connectivity” or “the team will learn how to implement power control software” or similar."
if ($instrPara -ne $null) {
    $instrPara.Range.Delete()
}

# --- 3. Move the rendered-page-break hint from "Schedule" to "Derived Requirements"
$schedulePara = Get-ParaByText $d "Schedule"
if ($schedulePara -ne $null) {
    $xmlSchedule = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="1647BBBF" w14:textId="77777777" w:rsidR="00DE5B15" w:rsidRDefault="00DE5B15" w:rsidP="00DE5B15">
<w:pPr><w:pStyle w:val="Heading1"/></w:pPr>
<w:r><w:t>Schedule</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
    $schedulePara.Range.InsertXML($xmlSchedule)
}

$derivedPara = Get-ParaByText $d "Derived Requirements"
if ($derivedPara -ne $null) {
    $xmlDerived = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="15E85400" w14:textId="77777777" w:rsidR="00DE5B15" w:rsidRDefault="00DE5B15" w:rsidP="00DE5B15">
<w:pPr><w:pStyle w:val="Heading1"/></w:pPr>
<w:r><w:lastRenderedPageBreak/><w:t>Derived Requirements</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
    $derivedPara.Range.InsertXML($xmlDerived)
}

Write-Output "done"
